# Registree stats backup on Sun 18 Apr 2021 21:32:10 SAST
#
# A new registree (Paijmans, Bronwyn Anne / Cowies Hill / 410E) is added to
# the MD410 and 410E attendance lists, and the "as of" timestamps plus the
# "Number of attendees" totals are refreshed on all affected sheets.

$wb = $excel.ActiveWorkbook

$newTimestamp = "18/04/2021 21:32"

# ---------------------------------------------------------------------
# Sheet 1: "MD410 Attendance" - full registree list (A:F, incl. District)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MD410 Attendance")

# Refresh title timestamp.
$ws1.Range("A1").Value = "MD410 Registrees as of $newTimestamp"

# Insert the new registree row, alphabetically before "Pantoleon" (row 144),
# pushing all subsequent rows down by one.
$ws1.Rows.Item(144).Insert()

# Copy formatting (style + row height) from the row that was just pushed
# down to row 145, then fill in the new values.
$ws1.Range("A145:F145").Copy()
$ws1.Range("A144:F144").PasteSpecial(-4122)
$ws1.Rows.Item(144).RowHeight = 25

$ws1.Cells.Item(144, 1).Value = "Paijmans"
$ws1.Cells.Item(144, 2).Value = "Bronwyn Anne"
$ws1.Cells.Item(144, 3).Value = "Cowies Hill"
$ws1.Cells.Item(144, 4).Value = "No"
$ws1.Cells.Item(144, 5).Value = "No"
$ws1.Cells.Item(144, 6).Value = "410E"

# The trailing summary rows shifted from 232/233 to 233/234; bump the
# attendee count (229 -> 230). The voter count line is unaffected.
$ws1.Range("A233").Value = "Number of attendees: 230"

# ---------------------------------------------------------------------
# Sheet 2: "410E Attendance" - district-specific list (A:E, no District)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("410E Attendance")

$ws2.Range("A1").Value = "410E Registrees as of $newTimestamp"

$ws2.Rows.Item(76).Insert()

$ws2.Range("A77:E77").Copy()
$ws2.Range("A76:E76").PasteSpecial(-4122)
$ws2.Rows.Item(76).RowHeight = 25

$ws2.Cells.Item(76, 1).Value = "Paijmans"
$ws2.Cells.Item(76, 2).Value = "Bronwyn Anne"
$ws2.Cells.Item(76, 3).Value = "Cowies Hill"
$ws2.Cells.Item(76, 4).Value = "No"
$ws2.Cells.Item(76, 5).Value = "No"

# Trailing summary rows shifted from 123/124 to 124/125; bump the
# attendee count (120 -> 121). The voter count line is unaffected.
$ws2.Range("A124").Value = "Number of attendees: 121"

# ---------------------------------------------------------------------
# Sheet 3: "410W Attendance" - only the "as of" timestamp changes.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("410W Attendance")
$ws3.Range("A1").Value = "410W Registrees as of $newTimestamp"

# ---------------------------------------------------------------------
# Sheet 4: "410E Voting" - only the "as of" timestamp changes.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("410E Voting")
$ws4.Range("A1").Value = "410E Voting details as of $newTimestamp"

# ---------------------------------------------------------------------
# Sheet 5: "410W Voting" - only the "as of" timestamp changes.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("410W Voting")
$ws5.Range("A1").Value = "410W Voting details as of $newTimestamp"
